# Delete slide 7 ("Ograniczenia") from the presentation.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$s.Delete()
